$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the npc_speed_gt_dut_start range value
$ws.Range("E3").Value = "[-10..10]kph"

# Update the npc_dist_gt_dut_end range value
$ws.Range("H3").Value = "[-5..15]m"

# Move the selection/active cell to H3 as seen in the saved workbook
$ws.Range("H3").Select()
